# Insert a new bold paragraph "Following are the features" right after
# the "Features " heading paragraph (and before the trailing empty
# paragraph that follows it), matching the heading's formatting.

$d = $word.ActiveDocument

# Locate the "Features" paragraph by scanning the paragraph collection
# (more robust than a hard-coded index).
$featuresPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -match "Features") {
        $featuresPara = $candidate
        break
    }
}

# Insert a brand-new paragraph right after it. InsertParagraphAfter()
# clones the paragraph/run formatting of the "Features" paragraph, so
# the new paragraph already picks up the bold run + "both" justification.
$featuresPara.Range.InsertParagraphAfter()

# Find the freshly inserted (still-empty) paragraph and give it its text.
$newPara = $featuresPara.Next()
$newPara.Range.Text = "Following are the features"
